$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 73158
$ws.Range("C2").Value = 5416.988499999999
$ws.Range("D2").Value = 67741.01150000001

$ws.Range("B3").Value = 69971
$ws.Range("C3").Value = 5210.142
$ws.Range("D3").Value = 64760.858

$ws.Range("B4").Value = 67233
$ws.Range("C4").Value = 5131.699999999999
$ws.Range("D4").Value = 62101.3

$ws.Range("B5").Value = 66774
$ws.Range("C5").Value = 5049.197999999999
$ws.Range("D5").Value = 61724.802

$ws.Range("B6").Value = 68427
$ws.Range("C6").Value = 5084.9925
$ws.Range("D6").Value = 63342.0075

$ws.Range("B7").Value = 72833
$ws.Range("C7").Value = 5380.1055
$ws.Range("D7").Value = 67452.89449999999

$ws.Range("B8").Value = 72858
$ws.Range("C8").Value = 6413.1515
$ws.Range("D8").Value = 66444.84849999999

$ws.Range("B9").Value = 84685
$ws.Range("C9").Value = 7207.8405
$ws.Range("D9").Value = 77477.15949999999

$ws.Range("B10").Value = 102552
$ws.Range("C10").Value = 9489.200000000001
$ws.Range("D10").Value = 93062.8

$ws.Range("B11").Value = 112257
$ws.Range("C11").Value = 14318.1885
$ws.Range("D11").Value = 97938.8115

$ws.Range("B12").Value = 116157
$ws.Range("C12").Value = 15340.192
$ws.Range("D12").Value = 100816.808

$ws.Range("B13").Value = 115181
$ws.Range("C13").Value = 15589.8715
$ws.Range("D13").Value = 99591.12850000001

$ws.Range("B14").Value = 118438
$ws.Range("C14").Value = 15359.3685
$ws.Range("D14").Value = 103078.6315

$ws.Range("B15").Value = 118968
$ws.Range("C15").Value = 15511.44
$ws.Range("D15").Value = 103456.56

$ws.Range("B16").Value = 121806
$ws.Range("C16").Value = 15666.4375
$ws.Range("D16").Value = 106139.5625

$ws.Range("B17").Value = 114966
$ws.Range("C17").Value = 16125.417
$ws.Range("D17").Value = 98840.583

$ws.Range("B18").Value = 15556
$ws.Range("C18").Value = 16406.873
$ws.Range("D18").Value = 0

$ws.Range("C19").Value = 16163.1645
$ws.Range("D19").Value = 91921.8355

$ws.Range("C20").Value = 15659.7175
$ws.Range("D20").Value = 91677.2825

$ws.Range("C21").Value = 13764.5865
$ws.Range("D21").Value = 78811.4135

$ws.Range("C22").Value = 11506.208
$ws.Range("D22").Value = 80527.792

$ws.Range("C23").Value = 8102.793999999999
$ws.Range("D23").Value = 59160.206

$ws.Range("C24").Value = 6308.112999999999
$ws.Range("D24").Value = 65087.887

$ws.Range("C25").Value = 5711.482
